$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values, entered so unique shared-strings are created in the
# required order (k, ㅇ, ㅁ, ㄴ, ㄹ, ㅎ): row 5 is written before row 4.
$ws.Range("A1").Value = "k"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = "ㅇ"
$ws.Range("A5").Value = "ㅁ"
$ws.Range("A4").Value = "ㄴ"
$ws.Range("A6").Value = "ㅇ"
$ws.Range("A7").Value = "ㄹ"
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 1
$ws.Range("A10").Value = "ㄴ"
$ws.Range("A11").Value = "ㄴ"
$ws.Range("A12").Value = "ㄴ"
$ws.Range("A13").Value = "ㄴ"

# Column B values (row 1..13)
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = "ㅎ"
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = "ㄹ"
$ws.Range("B10").Value = "ㅎ"
$ws.Range("B11").Value = "ㅁ"
$ws.Range("B12").Value = "ㅎ"
$ws.Range("B13").Value = "ㄴ"

$ws.Range("B13").Select()
